$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig14")

$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A112").Value = "Source: Short-Term Energy Outlook, February 2017."

$ws.Range("B86").Value = 488.64
$ws.Range("B87").Value = 480.18257143
$ws.Range("B88").Value = 495.97141527000002
$ws.Range("B89").Value = 499.62909999999999
$ws.Range("B90").Value = 509.1182
$ws.Range("B91").Value = 514.06349999999998
$ws.Range("B92").Value = 510.18770000000001
$ws.Range("B93").Value = 495.99779999999998
$ws.Range("B94").Value = 478.88389999999998
$ws.Range("B95").Value = 470.16770000000002
$ws.Range("B96").Value = 469.62369999999999
$ws.Range("B97").Value = 476.0222
$ws.Range("B98").Value = 471.7328
$ws.Range("B99").Value = 456.81009999999998
$ws.Range("B100").Value = 466.6053
$ws.Range("B101").Value = 475.5625
$ws.Range("B102").Value = 487.65910000000002
$ws.Range("B103").Value = 493.96969999999999
$ws.Range("B104").Value = 491.36380000000003
$ws.Range("B105").Value = 480.31380000000001
$ws.Range("B106").Value = 466.00119999999998
$ws.Range("B107").Value = 459.56400000000002
$ws.Range("B108").Value = 458.86869999999999
$ws.Range("B109").Value = 465.92720000000003
$ws.Range("B110").Value = 464.16449999999998
$ws.Range("B111").Value = 453.60359999999997

$ws.Range("A116").Value = 61
$ws.Range("A117").Value = 61
